# Realestate Update resale numbers 2025-02-05 21:44
# Appends a new data row (row 53) to the CityResaleNum sheet, mirroring the
# existing rows' layout: columns A-D are text, columns E-T are numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53

# Format the text columns (A-D) as Text first so that Excel doesn't
# auto-convert the date/time-looking strings (or the zero-padded week
# number) into numeric/date values. The format is reset back to Normal
# afterwards so the new row matches the plain styling of the other rows.
$textRange = $ws.Range("A${row}:D${row}")
$textRange.NumberFormat = "@"

$ws.Range("A${row}").Value = "2025-02-05"
$ws.Range("B${row}").Value = "21:44:17"
$ws.Range("C${row}").Value = "Wednesday"
$ws.Range("D${row}").Value = "05"

$textRange.Style = "Normal"

# Numeric columns (E-T)
$ws.Range("E${row}").Value = 125834
$ws.Range("F${row}").Value = 141749
$ws.Range("G${row}").Value = 167327
$ws.Range("H${row}").Value = 157967
$ws.Range("I${row}").Value = -1
$ws.Range("J${row}").Value = 142418
$ws.Range("K${row}").Value = -1
$ws.Range("L${row}").Value = -1
$ws.Range("M${row}").Value = 191102
$ws.Range("N${row}").Value = 115149
$ws.Range("O${row}").Value = 44748
$ws.Range("P${row}").Value = 28234
$ws.Range("Q${row}").Value = 63301
$ws.Range("R${row}").Value = -1
$ws.Range("S${row}").Value = 40168
$ws.Range("T${row}").Value = -1
